# Actualización automática 2025-10-07 17:30:09
# Applies the data updates for GUERRERO FAREZ FABIAN MAURICIO across the
# three worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (row 52: GUERRERO FAREZ FABIAN MAURICIO /
# TORO BLACIO MARIA DEL CISNE; row 56: count-of-54 summary row)
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D52").Value = 1408.32
$wsGrupo.Range("I52").Value = 127.8
$wsGrupo.Range("D56").Value = "3 de 54"
$wsGrupo.Range("I56").Value = "5 de 54"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (rows 53 & 54: GUERRERO FAREZ FABIAN MAURICIO /
# TORO BLACIO MARIA DEL CISNE "octubre" column; row 60: totals row)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F53").Value = 1979.47
$wsMensual.Range("F54").Value = 1979.47
$wsMensual.Range("F60").Value = 10832.04

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (row 3: 240X80 PORCELANATO; row 7:
# LAVABOS; row 14: TOTAL) - VENTA, POR CUMPLIR and CUMPLIMIENTO columns
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D3").Value = 2782.08
$wsCumpl.Range("E3").Value = 14887.0670988183
$wsCumpl.Range("F3").Value = 0.1574541195701553

$wsCumpl.Range("D7").Value = 548.1
$wsCumpl.Range("E7").Value = 338.611016287574
$wsCumpl.Range("F7").Value = 0.6181269770333413

$wsCumpl.Range("D14").Value = 8811.35
$wsCumpl.Range("E14").Value = 90205.15661190613
$wsCumpl.Range("F14").Value = 0.08898869796059325
